$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2024-05-15", "12:56:48", "Ascensor no sube",   "-", "-", "-", "-", "12:56:57"),
    @("2024-05-15", "12:57:05", "Fallo atornillador",  "-", "-", "-", "-", "12:57:23"),
    @("2024-05-15", "12:57:25", "Fallo atornillador",  "-", "-", "-", "-", "12:57:31"),
    @("2024-05-15", "12:57:26", "Fallo atornillador",  "-", "-", "-", "-", "12:57:32"),
    @("2024-05-15", "12:57:28", "Fallo atornillador",  "-", "-", "-", "-", "12:57:32")
)

$startRow = 70
$endRow = $startRow + $data.Count - 1

$textRange = $ws.Range("A$startRow`:H$endRow")
$textRange.NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $vals[$c]
    }
}
